$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are written as Text, matching the source
# data which uses locale-formatted numbers (e.g. "29.492.95") that must not
# be reinterpreted as numeric values by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.492.95'
$ws.Range("E2").Value = '  +0.91%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.919.92'
$ws.Range("E3").Value = '  +1.55%  '

$ws.Range("E4").Value = '  +0.71%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.38'

$ws.Range("E6").Value = '  +0.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4830'
$ws.Range("E7").Value = '  +2.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4079'
$ws.Range("E8").Value = '  +1.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08225'
$ws.Range("E9").Value = '  +2.66%  '

$ws.Range("E10").Value = '  +2.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.57'
$ws.Range("E11").Value = '  +3.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.934.30'
$ws.Range("E12").Value = '  -1.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.053'
$ws.Range("E13").Value = '  +2.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.237'
$ws.Range("E14").Value = '  +2.83%  '

$ws.Range("E15").Value = '  +2.42%  '

$ws.Range("E16").Value = '  +2.71%  '

$ws.Range("E17").Value = '  +0.66%  '

$ws.Range("E18").Value = '  +1.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.82'
$ws.Range("E19").Value = '  +1.93%  '

$ws.Range("E20").Value = '  +0.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.502.81'
$ws.Range("E21").Value = '  +0.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.648'
$ws.Range("E22").Value = '  +2.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.77'
$ws.Range("E23").Value = '  +0.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.195'
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.151.70'
$ws.Range("E25").Value = '  -1.83%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.618'
$ws.Range("E26").Value = '  +11.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.65'
$ws.Range("E27").Value = '  +0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.14'
$ws.Range("E28").Value = '  +2.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.119'
$ws.Range("E29").Value = '  +1.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.52'
$ws.Range("E30").Value = '  +2.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.025'
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("E32").Value = '  +1.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.545'
$ws.Range("E33").Value = '  +3.65%  '

$ws.Range("E34").Value = '  +0.96%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.383'
$ws.Range("E35").Value = '  +0.30%  '

$ws.Range("E36").Value = '  +1.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06141'
$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("E38").Value = '  +0.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5985'
$ws.Range("E39").Value = '  +2.88%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.86'
$ws.Range("E40").Value = '  +8.39%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.033'
$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1857'
$ws.Range("E42").Value = '  +1.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.417'
$ws.Range("E43").Value = '  +2.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.280'
$ws.Range("E44").Value = '  -0.24%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.47'
$ws.Range("E45").Value = '  +2.00%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07603'
$ws.Range("E46").Value = '  -1.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5583'
$ws.Range("E47").Value = '  +2.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.959'
$ws.Range("E48").Value = '  +2.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.81'
$ws.Range("E49").Value = '  +4.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.427'
$ws.Range("E50").Value = '  +4.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.83'
$ws.Range("E51").Value = '  +2.41%  '
